# CodeGenerators/out/parsed.xlsx - "Other" (.1 id) rows removed from the
# questionnaire sheet; the PLT value that had been sitting on each removed
# "Other" row is merged up onto the row directly above it before the row
# is deleted.
#
# Rows below are 1-indexed worksheet rows (row 1 is the header ROW/DATATYPE/
# ID/TEXT/PLT row), matching the *original* (before-edit) layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Merge the PLT ("E") value of each doomed "Other" row onto the row
#    immediately above it (its parent question), while row numbers still
#    match the original layout.
$ws.Cells.Item(5, 5).Value  = 355   # CQ3   <- CQ3.1
$ws.Cells.Item(10, 5).Value = 162   # CQ4b  <- CQ4b.1
$ws.Cells.Item(20, 5).Value = 154   # CQ7a  <- CQ7a.1
$ws.Cells.Item(22, 5).Value = 162   # CQ5   <- CQ5.1
$ws.Cells.Item(29, 5).Value = 162   # CQ7   <- CQ7.1
$ws.Cells.Item(63, 5).Value = 162   # CQ27  <- CQ27.1
$ws.Cells.Item(67, 5).Value = 157   # CQ29A <- CQ29A.1
# (row 4 -> CQ2.1 and row 38 -> CQ11.1 merge a PLT of 0 onto a PLT of 0 -
#  no value change needed there.)

# 2) Delete the "Other" rows themselves, bottom-to-top so the remaining,
#    not-yet-deleted row numbers stay valid.
$rowsToDelete = @(68, 64, 38, 30, 23, 21, 11, 6, 4)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# 3) Column A ("ROW") is a literal sequence number, not a formula, so it
#    does not renumber itself when rows are deleted - rewrite it to stay
#    contiguous (1..85) for the surviving data rows (2..86).
for ($i = 2; $i -le 86; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

Write-Output "done"
